$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035343593985242
$ws.Range("D2").Value = 1.044958576889714
$ws.Range("E2").Value = 1.043646377795473
$ws.Range("F2").Value = 1.052708748422648
$ws.Range("I2").Value = 1.026260528590535
$ws.Range("J2").Value = 1.040457885307983
$ws.Range("K2").Value = 1.047728163742146
$ws.Range("L2").Value = 1.046419656740346
$ws.Range("M2").Value = 1.055456732356759
$ws.Range("N2").Value = 1.017438604218481
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036656968303472
$ws.Range("D3").Value = 1.046249757366403
$ws.Range("E3").Value = 1.044870185010141
$ws.Range("F3").Value = 1.054167554089118
$ws.Range("I3").Value = 1.026437532419498
$ws.Range("J3").Value = 1.041412932363572
$ws.Range("K3").Value = 1.048829437406635
$ws.Range("L3").Value = 1.047453462874419
$ws.Range("M3").Value = 1.0567267815579
$ws.Range("N3").Value = 1.017753466389779
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037500774956719
$ws.Range("D4").Value = 1.047077657355115
$ws.Range("E4").Value = 1.045654924263664
$ws.Range("F4").Value = 1.055098944051483
$ws.Range("I4").Value = 1.026542394296678
$ws.Range("J4").Value = 1.042024454768805
$ws.Range("K4").Value = 1.04953392684855
$ws.Range("L4").Value = 1.048114715372643
$ws.Range("M4").Value = 1.057535549989291
$ws.Range("N4").Value = 1.017955063093808
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037854081949189
$ws.Range("D5").Value = 1.04742390836385
$ws.Range("E5").Value = 1.045983133442264
$ws.Range("F5").Value = 1.055487517790941
$ws.Range("I5").Value = 1.026584167974395
$ws.Range("J5").Value = 1.042280005883959
$ws.Range("K5").Value = 1.049828168557979
$ws.Range("L5").Value = 1.048390879693913
$ws.Range("M5").Value = 1.05787245617933
$ws.Range("N5").Value = 1.018039306237648
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037913320361235
$ws.Range("D6").Value = 1.047481940449837
$ws.Range("E6").Value = 1.046038142233479
$ws.Range("F6").Value = 1.055552586772024
$ws.Range("I6").Value = 1.026591046616384
$ws.Range("J6").Value = 1.042322824510152
$ws.Range("K6").Value = 1.049877460545901
$ws.Range("L6").Value = 1.048437142213823
$ws.Range("M6").Value = 1.057928843014198
$ws.Range("N6").Value = 1.018053421355035
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037505501457944
$ws.Range("D7").Value = 1.047082291021597
$ws.Range("E7").Value = 1.045659316451188
$ws.Range("F7").Value = 1.055104147881976
$ws.Range("I7").Value = 1.026542961549701
$ws.Range("J7").Value = 1.042027875461504
$ws.Range("K7").Value = 1.04953786606669
$ws.Range("L7").Value = 1.048118412650367
$ws.Range("M7").Value = 1.057540063899285
$ws.Range("N7").Value = 1.017956190745712
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035788713927501
$ws.Range("D8").Value = 1.045396517182373
$ws.Range("E8").Value = 1.044061458427446
$ws.Range("F8").Value = 1.053204372569413
$ws.Range("I8").Value = 1.026322353747458
$ws.Range("J8").Value = 1.040781992898821
$ws.Range("K8").Value = 1.048102032847066
$ws.Range("L8").Value = 1.046770637548922
$ws.Range("M8").Value = 1.055888664442549
$ws.Range("N8").Value = 1.017545459145822
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032716567673166
$ws.Range("D9").Value = 1.042367116990732
$ws.Range("E9").Value = 1.041190351809656
$ws.Range("F9").Value = 1.04975945398698
$ws.Range("I9").Value = 1.025859266304263
$ws.Range("J9").Value = 1.038536514901501
$ws.Range("K9").Value = 1.045509063356386
$ws.Range("L9").Value = 1.044336090438673
$ws.Range("M9").Value = 1.052877790425674
$ws.Range("N9").Value = 1.016805101568613
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03063581767105
$ws.Range("D10").Value = 1.040306755104612
$ws.Range("E10").Value = 1.039237860079458
$ws.Range("F10").Value = 1.047395759937005
$ws.Range("I10").Value = 1.02550009595894
$ws.Range("J10").Value = 1.037004930215756
$ws.Range("K10").Value = 1.043737036393813
$ws.Range("L10").Value = 1.042671929910005
$ws.Range("M10").Value = 1.050801130870269
$ws.Range("N10").Value = 1.016300064518186
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029726824456937
$ws.Range("D11").Value = 1.039404639161595
$ws.Range("E11").Value = 1.038383025616001
$ws.Range("F11").Value = 1.046355934814718
$ws.Range("I11").Value = 1.025332484181779
$ws.Range("J11").Value = 1.036333304937137
$ws.Range("K11").Value = 1.042959166872508
$ws.Range("L11").Value = 1.041941317322026
$ws.Range("M11").Value = 1.049885067031217
$ws.Range("N11").Value = 1.016078584211453
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029387956053361
$ws.Range("D12").Value = 1.039068029524192
$ws.Range("E12").Value = 1.038064066074196
$ws.Range("F12").Value = 1.045967206530262
$ws.Range("I12").Value = 1.025268397581603
$ws.Range("J12").Value = 1.036082544861732
$ws.Range("K12").Value = 1.04266861840073
$ws.Range("L12").Value = 1.041668406238501
$ws.Range("M12").Value = 1.049542232205205
$ws.Range("N12").Value = 1.015995889706036
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029460700432865
$ws.Range("D13").Value = 1.039140302819453
$ws.Range("E13").Value = 1.038132549371817
$ws.Range("F13").Value = 1.046050703412368
$ws.Range("I13").Value = 1.025282227304522
$ws.Range("J13").Value = 1.036136392402502
$ws.Range("K13").Value = 1.042731015447395
$ws.Range("L13").Value = 1.041727016178076
$ws.Range("M13").Value = 1.049615888234321
$ws.Range("N13").Value = 1.016013647387536
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029698838678197
$ws.Range("D14").Value = 1.039376846142725
$ws.Range("E14").Value = 1.038356689762217
$ws.Range("F14").Value = 1.046323853452373
$ws.Range("I14").Value = 1.025327224129427
$ws.Range("J14").Value = 1.036312603440524
$ws.Range("K14").Value = 1.042935183101673
$ws.Range("L14").Value = 1.041918789739954
$ws.Range("M14").Value = 1.049856780794706
$ws.Range("N14").Value = 1.016071757406246
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029845400079313
$ws.Range("D15").Value = 1.039522385585611
$ws.Range("E15").Value = 1.038494599011156
$ws.Range("F15").Value = 1.046491818976285
$ws.Range("I15").Value = 1.02535470553927
$ws.Range("J15").Value = 1.036421001495326
$ws.Range("K15").Value = 1.043060763008693
$ws.Range("L15").Value = 1.042036744419912
$ws.Range("M15").Value = 1.050004861282476
$ws.Range("N15").Value = 1.016107504134339
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030695973696437
$ws.Range("D16").Value = 1.040366413383013
$ws.Range("E16").Value = 1.03929439264582
$ws.Range("F16").Value = 1.047464422493807
$ws.Range("I16").Value = 1.025510964111279
$ws.Range("J16").Value = 1.037049324233116
$ws.Range("K16").Value = 1.043788436223339
$ws.Range("L16").Value = 1.042720205161082
$ws.Range("M16").Value = 1.050861568884182
$ws.Range("N16").Value = 1.016314703957202
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031227353890222
$ws.Range("D17").Value = 1.040893163436641
$ws.Range("E17").Value = 1.039793550264466
$ws.Range("F17").Value = 1.048070113119884
$ws.Range("I17").Value = 1.025605736410267
$ws.Range("J17").Value = 1.037441180867676
$ws.Range("K17").Value = 1.044242039211078
$ws.Range("L17").Value = 1.043146223077529
$ws.Range("M17").Value = 1.05139442157327
$ws.Range("N17").Value = 1.016443921691499
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031536527328094
$ws.Range("D18").Value = 1.041199448148912
$ws.Range("E18").Value = 1.040083795951815
$ws.Range("F18").Value = 1.048421829256416
$ws.Range("I18").Value = 1.025659849976952
$ws.Range("J18").Value = 1.037668931429029
$ws.Range("K18").Value = 1.044505600060295
$ws.Range("L18").Value = 1.043393746993889
$ws.Range("M18").Value = 1.051703601624537
$ws.Range("N18").Value = 1.01651902292214
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031641817297214
$ws.Range("D19").Value = 1.041303721318654
$ws.Range("E19").Value = 1.040182609541116
$ws.Range("F19").Value = 1.048541489672364
$ws.Range("I19").Value = 1.025678103958108
$ws.Range("J19").Value = 1.037746451231121
$ws.Range("K19").Value = 1.044595295548349
$ws.Range("L19").Value = 1.043477983245507
$ws.Range("M19").Value = 1.051808749629381
$ws.Range("N19").Value = 1.016544585025336
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031170421806903
$ws.Range("D20").Value = 1.040836747592241
$ws.Range("E20").Value = 1.039740089109044
$ws.Range("F20").Value = 1.04800529121179
$ws.Range("I20").Value = 1.025595688877107
$ws.Range("J20").Value = 1.037399222575712
$ws.Range("K20").Value = 1.044193477387365
$ws.Range("L20").Value = 1.043100615381648
$ws.Range("M20").Value = 1.051337419741181
$ws.Range("N20").Value = 1.016430085753343
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029628746991234
$ws.Range("D21").Value = 1.039307232307357
$ws.Range("E21").Value = 1.038290725848992
$ws.Range("F21").Value = 1.046243486647163
$ws.Range("I21").Value = 1.025314024252197
$ws.Range("J21").Value = 1.036260749423229
$ws.Range("K21").Value = 1.042875105553346
$ws.Range("L21").Value = 1.041862359610733
$ws.Range("M21").Value = 1.049785915121437
$ws.Range("N21").Value = 1.016054657293814
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028652318305435
$ws.Range("D22").Value = 1.038336736157941
$ws.Range("E22").Value = 1.037371132395184
$ws.Range("F22").Value = 1.0451213383452
$ws.Range("I22").Value = 1.025126346900654
$ws.Range("J22").Value = 1.035537480034655
$ws.Range("K22").Value = 1.0420368466036
$ws.Range("L22").Value = 1.041074959903986
$ws.Range("M22").Value = 1.048795546362824
$ws.Range("N22").Value = 1.015816137228439
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029170624411293
$ws.Range("D23").Value = 1.038852060597388
$ws.Range("E23").Value = 1.03785942366503
$ws.Range("F23").Value = 1.045717591427861
$ws.Range("I23").Value = 1.025226845679104
$ws.Range("J23").Value = 1.03592161371075
$ws.Range("K23").Value = 1.04248211816196
$ws.Range("L23").Value = 1.041493223351978
$ws.Range("M23").Value = 1.04932198202144
$ws.Range("N23").Value = 1.015942818030252
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03119614934828
$ws.Range("D24").Value = 1.040862242448368
$ws.Range("E24").Value = 1.039764248697977
$ws.Range("F24").Value = 1.048034586300721
$ws.Range("I24").Value = 1.025600232525836
$ws.Range("J24").Value = 1.037418184232495
$ws.Range("K24").Value = 1.044215423530136
$ws.Range("L24").Value = 1.04312122651758
$ws.Range("M24").Value = 1.05136318143188
$ws.Range("N24").Value = 1.016436338450077
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033516462787259
$ws.Range("D25").Value = 1.043157382013318
$ws.Range("E25").Value = 1.041939287047277
$ws.Range("F25").Value = 1.050661741442136
$ws.Range("I25").Value = 1.025987839603171
$ws.Range("J25").Value = 1.03912305114503
$ws.Range("K25").Value = 1.046186967246703
$ws.Range("L25").Value = 1.044972646474046
$ws.Range("M25").Value = 1.05366828547454
$ws.Range("N25").Value = 1.016998498887812
